$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Tnfrsf19"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 41.528285
$ws.Range("H2").Value = 124.584855
$ws.Range("I2").Value = 0.137866712381124
$ws.Range("J2").Value = 0.145939792231724
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04225366666666667
$ws.Range("N2").Value = 0.126761
$ws.Range("O2").Value = 0.03359366718256408
$ws.Range("P2").Value = 0.0495580819187537
$ws.Range("Q2").Value = 1.754722311628333
$ws.Range("R2").Value = 15.792500804655
$ws.Range("S2").Value = 0.004631448451285766
$ws.Range("T2").Value = 0.007232496178625671

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Tnfrsf19"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 41.528285
$ws.Range("H3").Value = 124.584855
$ws.Range("I3").Value = 0.137866712381124
$ws.Range("J3").Value = 0.145939792231724
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.215533
$ws.Range("N3").Value = 2.431066
$ws.Range("O3").Value = 0.9664063328174359
$ws.Range("P3").Value = 0.9504419180812462
$ws.Range("Q3").Value = 50.479000850905
$ws.Range("R3").Value = 302.87400510543
$ws.Range("S3").Value = 0.1332352639298382
$ws.Range("T3").Value = 0.1387072960530983

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Tnfrsf19"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 72.02213166666667
$ws.Range("H4").Value = 216.066395
$ws.Range("I4").Value = 0.2391010009578718
$ws.Range("J4").Value = 0.2531020708300187
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04225366666666667
$ws.Range("N4").Value = 0.126761
$ws.Range("O4").Value = 0.03359366718256408
$ws.Range("P4").Value = 0.0495580819187537
$ws.Range("Q4").Value = 3.043199144066111
$ws.Range("R4").Value = 27.388792296595
$ws.Range("S4").Value = 0.008032279449196681
$ws.Range("T4").Value = 0.01254325316000027

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Tnfrsf19"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 72.02213166666667
$ws.Range("H5").Value = 216.066395
$ws.Range("I5").Value = 0.2391010009578718
$ws.Range("J5").Value = 0.2531020708300187
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.215533
$ws.Range("N5").Value = 2.431066
$ws.Range("O5").Value = 0.9664063328174359
$ws.Range("P5").Value = 0.9504419180812462
$ws.Range("Q5").Value = 87.54527777117833
$ws.Range("R5").Value = 525.27166662707
$ws.Range("S5").Value = 0.2310687215086751
$ws.Range("T5").Value = 0.2405588176700184

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Tnfrsf19"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 65.334877
$ws.Range("H6").Value = 196.004631
$ws.Range("I6").Value = 0.2169004738773853
$ws.Range("J6").Value = 0.2296015444621718
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04225366666666667
$ws.Range("N6").Value = 0.126761
$ws.Range("O6").Value = 0.03359366718256408
$ws.Range("P6").Value = 0.0495580819187537
$ws.Range("Q6").Value = 2.760638114465667
$ws.Range("R6").Value = 24.84574303019101
$ws.Range("S6").Value = 0.007286482331177316
$ws.Range("T6").Value = 0.01137861214912868

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Tnfrsf19"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 65.334877
$ws.Range("H7").Value = 196.004631
$ws.Range("I7").Value = 0.2169004738773853
$ws.Range("J7").Value = 0.2296015444621718
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.215533
$ws.Range("N7").Value = 2.431066
$ws.Range("O7").Value = 0.9664063328174359
$ws.Range("P7").Value = 0.9504419180812462
$ws.Range("Q7").Value = 79.416699044441
$ws.Range("R7").Value = 476.500194266646
$ws.Range("S7").Value = 0.209613991546208
$ws.Range("T7").Value = 0.2182229323130431

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Tnfrsf19"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 72.34659833333335
$ws.Range("H8").Value = 217.039795
$ws.Range("I8").Value = 0.2401781740848285
$ws.Range("J8").Value = 0.2542423201304522
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04225366666666667
$ws.Range("N8").Value = 0.126761
$ws.Range("O8").Value = 0.03359366718256408
$ws.Range("P8").Value = 0.0495580819187537
$ws.Range("Q8").Value = 3.05690905044389
$ws.Range("R8").Value = 27.51218145399501
$ws.Range("S8").Value = 0.008068465644721665
$ws.Range("T8").Value = 0.01259976172823895

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Tnfrsf19"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 72.34659833333335
$ws.Range("H9").Value = 217.039795
$ws.Range("I9").Value = 0.2401781740848285
$ws.Range("J9").Value = 0.2542423201304522
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.215533
$ws.Range("N9").Value = 2.431066
$ws.Range("O9").Value = 0.9664063328174359
$ws.Range("P9").Value = 0.9504419180812462
$ws.Range("Q9").Value = 87.93967771191168
$ws.Range("R9").Value = 527.63806627147
$ws.Range("S9").Value = 0.2321097084401068
$ws.Range("T9").Value = 0.2416425584022132

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rtn4"
$ws.Range("C10").Value = "Tnfrsf19"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.988644
$ws.Range("H10").Value = 99.977288
$ws.Range("I10").Value = 0.1659536386987904
$ws.Range("J10").Value = 0.1171142723456333
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.04225366666666667
$ws.Range("N10").Value = 0.126761
$ws.Range("O10").Value = 0.03359366718256408
$ws.Range("P10").Value = 0.0495580819187537
$ws.Range("Q10").Value = 2.112203500694667
$ws.Range("R10").Value = 12.673221004168
$ws.Range("S10").Value = 0.005574991306182653
$ws.Range("T10").Value = 0.005803958702760124

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Rtn4"
$ws.Range("C11").Value = "Tnfrsf19"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 49.988644
$ws.Range("H11").Value = 99.977288
$ws.Range("I11").Value = 0.1659536386987904
$ws.Range("J11").Value = 0.1171142723456333
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.215533
$ws.Range("N11").Value = 2.431066
$ws.Range("O11").Value = 0.9664063328174359
$ws.Range("P11").Value = 0.9504419180812462
$ws.Range("Q11").Value = 60.762846407252
$ws.Range("R11").Value = 243.051385629008
$ws.Range("S11").Value = 0.1603786473926078
$ws.Range("T11").Value = 0.1113103136428731
